$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 94; this shifts existing rows 94-132 down to 95-133
$ws.Rows(94).Insert()

# Populate the newly inserted row 94 with the new record
$ws.Cells.Item(94, 1).Value = 7
$ws.Cells.Item(94, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(94, 3).Value = "Ñuble"
$ws.Cells.Item(94, 4).Value = 44468
$ws.Cells.Item(94, 5).Value = 16
$ws.Cells.Item(94, 6).Value = 100112032
$ws.Cells.Item(94, 7).Value = "Zapallo italiano"
$ws.Cells.Item(94, 8).Value = "Sin especificar"
$ws.Cells.Item(94, 9).Value = "Primera"
$ws.Cells.Item(94, 10).Value = 120
$ws.Cells.Item(94, 11).Value = 14000
$ws.Cells.Item(94, 12).Value = 15000
$ws.Cells.Item(94, 13).Value = 14500
$ws.Cells.Item(94, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(94, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(94, 16).Value = 290
$ws.Cells.Item(94, 17).Value = 50
$ws.Cells.Item(94, 18).Value = "Hortaliza"
